$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a "last changed" date-serial for every data
# row. All existing rows (2..last) are being bumped by one day:
# 46074 (2026-02-21) -> 46075 (2026-02-22).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 46075
}
